$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = 1430
$ws.Range("D4").Value = 1600
$ws.Range("F4").Value = 1950
$ws.Range("G4").Value = 2100

$ws.Range("C5").Value = 1580
$ws.Range("D5").Value = 1580
$ws.Range("F5").Value = 2100
$ws.Range("G5").Value = 2050

$ws.Range("B6").Value = 1350
$ws.Range("D6").Value = 1680
$ws.Range("E6").Value = 900
$ws.Range("G6").Value = 2120

$ws.Range("C7").Value = 1450
$ws.Range("D7").Value = 1400
$ws.Range("G7").Value = 850

$ws.Range("C8").Value = 1550
$ws.Range("D8").Value = 1270
$ws.Range("E8").Value = 1870
$ws.Range("G8").Value = 750

$ws.Range("B9").Value = 1350
$ws.Range("D9").Value = 1250
$ws.Range("E9").Value = 1830
$ws.Range("G9").Value = 750

$ws.Range("C14").Select()
